$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.209.20"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.65%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.795.63"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +7.41%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "428.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +9.53%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +14.01%  "

# Row 7
$ws.Range("E7").Value = "  +5.59%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.01%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.747"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +10.71%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.156"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.34%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000323"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.02%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.42"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +12.88%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.64"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +16.96%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.394.23"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.89%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.03"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +15.32%  "

# Row 16
$ws.Range("E16").Value = "  +1.30%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.784.82"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.04%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +8.70%  "

# Row 19
$ws.Range("E19").Value = "  +12.42%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "66.387.59"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.03%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "411.82"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.34%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.34"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +10.75%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.29"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +14.04%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.73"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.75%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "37.17"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +10.98%  "

# Row 26
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.30"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +11.57%  "

# Row 27
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.75"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +45.87%  "

# Row 28
$ws.Range("E28").Value = "  +14.53%  "

# Row 29
$ws.Range("E29").Value = "  -0.73%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "13.90"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +18.98%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "703.34"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.48%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.129"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +18.24%  "

# Row 33
$ws.Range("E33").Value = "  +5.19%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "40.65"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +11.97%  "

# Row 35
$ws.Range("B35").Value = "Dai"
$ws.Range("C35").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.06%  "

# Row 36
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.75"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +41.60%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.152"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.08%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "56.48"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.10%  "

# Row 39
$ws.Range("E39").Value = "  +10.22%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.61"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +50.20%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0₃0684"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.90%  "

# Row 42
$ws.Range("B42").Value = "ThetaToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.86"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +8.02%  "

# Row 43
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.141"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +8.99%  "

# Row 44
$ws.Range("E44").Value = "  +0.48%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.37"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +10.80%  "

# Row 46
$ws.Range("E46").Value = "  +18.43%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.19"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.19%  "

# Row 48
$ws.Range("E48").Value = "  +6.96%  "

# Row 49
$ws.Range("E49").Value = "  +8.21%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "142.88"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.57%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.83"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.59%  "
